# Weekly update: insert a new daily price record for "Apio" at row 603,
# pushing the existing records (rows 603:691) down by one row to (604:692).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 603 (shifts 603:691 -> 604:692, and the sheet
# dimension grows from A1:R691 to A1:R692 automatically).
$ws.Rows.Item(603).Insert()

# Populate the newly inserted row 603 with the new weekly record.
$ws.Cells.Item(603, 1).Value = 8
$ws.Cells.Item(603, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(603, 3).Value = "Coquimbo"
$ws.Cells.Item(603, 4).Value = 45154
$ws.Cells.Item(603, 5).Value = 4
$ws.Cells.Item(603, 6).Value = 100112017
$ws.Cells.Item(603, 7).Value = "Apio"
$ws.Cells.Item(603, 8).Value = "Americana (o)"
$ws.Cells.Item(603, 9).Value = "Primera"
$ws.Cells.Item(603, 10).Value = 900
$ws.Cells.Item(603, 11).Value = 7000
$ws.Cells.Item(603, 12).Value = 8000
$ws.Cells.Item(603, 13).Value = 7500
$ws.Cells.Item(603, 14).Value = "`$/docena de matas"
$ws.Cells.Item(603, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(603, 16).Value = 1250
$ws.Cells.Item(603, 17).Value = 6
$ws.Cells.Item(603, 18).Value = "Hortaliza"
